# Update cryptocurrency Price (D) and Volume(1h) (E) columns with
# freshly scraped coinranking.com figures (GitHub Actions refresh run).
#
# Column D values are stored as plain text in the workbook (some of them
# use a dotted/European-style grouping like '67.248.62' that is not a
# legal number, so the whole column is kept as text for consistency).
# Assigning a string through Range.Value normally auto-converts a
# 'clean looking' numeric string (e.g. '0.999') into a real number, so
# each Price cell is briefly forced to Text format, written, and then
# restored to the workbook's default ('Normal') style - this keeps the
# value a text string without leaving a custom number format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.248.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.69%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.666.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.15%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.668.54"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.02%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  -0.80%  "

$ws.Range("E10").Value = "  -5.64%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.17"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.61%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000239"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.07%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.278.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.660.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.44%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.250.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.23%  "

$ws.Range("E20").Value = "  -3.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "492.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.95%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.14"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.81%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.717"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.33%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000136"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.24%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.20%  "

$ws.Range("E30").Value = "  -2.21%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.12%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.802.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.37%  "

$ws.Range("E35").Value = "  -5.90%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.601.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.16%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.989"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.11%  "

$ws.Range("E40").Value = "  -6.95%  "

$ws.Range("E41").Value = "  -4.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "432.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.76%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.47%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.60%  "

$ws.Range("E47").Value = "  -0.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.744.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0345"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.25%  "

